# ----------------------------------------------------------------------
# Applies the "raven.docx" commit:
#   1. Appends a red-colored annotation to the first paragraph:
#        "This is a Microsoft word document.  (This is a change – Version for main branch)"
#      where the trailing annotation is split into three runs (as Word's
#      autocorrect / incremental typing would naturally produce) all
#      colored red (FF0000), and the original sentence keeps two
#      trailing spaces appended to it in its own (unformatted) run.
#   2. Removes the trailing "ank God almighty, we are free at last."
#      paragraph (the tail half of a spell-checked word broken across
#      paragraphs) that followed "Shall be lifted-nevermore!".
#   3. Removes a set of unused styles (Heading 2/4 [+ their linked
#      Char styles], Hyperlink, apple-converted-space, audio-tool,
#      subscribe, subscribe-more-info, generic-title and
#      podcast-tools__subscribe-links) that were cleaned out of
#      styles.xml in the same commit.
# ----------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. First paragraph: split "This is a Microsoft word document." -------
$p1 = $d.Paragraphs(1).Range

# 1a. Two trailing spaces appended to the existing (plain) run.
$pos = $p1.End - 1
$rng = $d.Range($pos, $pos)
$rng.InsertAfter("  ")
$pos = $p1.End - 1

# 1b. Red run: "(This is a change – Ve"
$rng = $d.Range($pos, $pos)
$rng.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$newPos = $p1.End - 1
$d.Range($pos, $newPos).Font.Color = 255
$pos = $newPos

# 1c. Red run: "rsion for main branch"
$rng = $d.Range($pos, $pos)
$rng.InsertAfter("rsion for main branch")
$newPos = $p1.End - 1
$d.Range($pos, $newPos).Font.Color = 255
$pos = $newPos

# 1d. Red run: ")"
$rng = $d.Range($pos, $pos)
$rng.InsertAfter(")")
$newPos = $p1.End - 1
$d.Range($pos, $newPos).Font.Color = 255
$pos = $newPos

# --- 2. Drop the trailing "ank God almighty, we are free at last." para ---
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)
if ($lastPara.Range.Text -like "*God almighty, we are free at last.*") {
    $lastPara.Range.Delete()
}

# --- 3. Remove the now-unused styles ---------------------------------------
$stylesToRemove = @(
    "Heading 2",
    "Heading 4",
    "Heading 2 Char",
    "Heading 4 Char",
    "Hyperlink",
    "apple-converted-space",
    "audio-tool",
    "subscribe",
    "subscribe-more-info",
    "generic-title",
    "podcast-tools__subscribe-links"
)

$i = $d.Styles.Count
while ($i -ge 1) {
    $s = $d.Styles.Item($i)
    if ($stylesToRemove -contains $s.NameLocal) {
        $s.Delete()
    }
    $i = $i - 1
}

Write-Output "edit complete"
